$wb = $excel.ActiveWorkbook

# --- Hoja1!A1: update the two "conversión del día" rate lines ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 8.61 = 36022.38 pesos`n✅ 36022.38 pesos = 8.59 = 946.55 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$ws1.Range("A1").Value = $newText

# --- tasas sheet: update the N10/O10 and N12/O12 rate values ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 116.15
$ws2.Range("O10").Value = 4184
$ws2.Range("N12").Value = 4194
$ws2.Range("O12").Value = 110.205
